$d = $word.ActiveDocument

# Locate the start of the paragraph containing the closing quote text "apache2ctl -S"
# (unique first hit) and then, searching *after* that point, the "autoload-psr4.php"
# sub-title paragraph (the same text also appears earlier, in the table of contents,
# so the second search must start after the first match). Replace that whole span
# (the quote paragraph through the page break, "Modifications" title, the blank line,
# and the "autoload-psr4.php" sub-title) with the corrected OOXML in one shot via
# Range.InsertXML -- this merges the quote paragraph with the page-break paragraph and
# strips the now-redundant w:lang markup, matching the authored edit.

$startRng = $d.Content
$null = $startRng.Find.Execute("apache2ctl -S", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $startRng.Paragraphs(1)

$endRng = $d.Range($startPara.Range.End, $d.Content.End)
$null = $endRng.Find.Execute("autoload-psr4.php", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPara = $endRng.Paragraphs(1)

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<w:p w14:paraId="7D1A1AB9" w14:textId="0FBAC8DA" w:rsidR="00A14FCF" w:rsidRPr="005E4A1C" w:rsidRDefault="00B52A41"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="005E4A1C"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t>« </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005E4A1C"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidR="00A14FCF" w:rsidRPr="005E4A1C"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t>udo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00A14FCF" w:rsidRPr="005E4A1C"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> apache2ctl -S</w:t></w:r><w:r w:rsidRPr="005E4A1C"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t> »</w:t></w:r><w:bookmarkStart w:id="10" w:name="_Toc131579627"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p w14:paraId="00000047" w14:textId="674098CA" w:rsidR="00802F63" w:rsidRPr="00473AC2" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Titre"/><w:rPr><w:rFonts w:cs="Arial"/></w:rPr></w:pPr><w:r w:rsidRPr="00473AC2"><w:rPr><w:rFonts w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t>Modifications</w:t></w:r><w:bookmarkEnd w:id="10"/></w:p><w:p w14:paraId="00000048" w14:textId="77777777" w:rsidR="00802F63" w:rsidRPr="00473AC2" w:rsidRDefault="00802F63"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p w14:paraId="00000049" w14:textId="77777777" w:rsidR="00802F63" w:rsidRPr="00473AC2" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Sous-titre"/></w:pPr><w:bookmarkStart w:id="11" w:name="_Toc131579628"/><w:r w:rsidRPr="00473AC2"><w:t>autoload-psr4.php</w:t></w:r><w:bookmarkEnd w:id="11"/></w:p>'

$target.InsertXML($xml)

Write-Output "done"
